$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")
$ws.Range("A99").Value = "WUR2013Treatment60WM"
$ws.Range("A106").Value = "WUR2014Treatment60WM"
$ws.Range("A119").Value = "WUR2014Treatment02WM"
$ws.Range("A111").Value = "WUR2013Treatment02WM"
$ws.Range("A125").Value = "WUR2013Treatment62WM"
$ws.Range("A137").Value = "WUR2014Treatment62WM"
